$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------------------
# Locate the existing "TASK SCHEDULER Tips:" paragraph. All new content is
# spliced in immediately before its first run (new paragraphs 1-2) and
# immediately after its closing </w:p> (new paragraphs 3-4).
# ---------------------------------------------------------------------------
$findR = $d.Content
$found = $findR.Find.Execute("TASK SCHEDULER Tips:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate 'TASK SCHEDULER Tips:' anchor text"
}

$taskPara = $findR.Paragraphs.Item(1)
$taskStart = $taskPara.Range.Start
$taskEnd = $taskPara.Range.End

# ---------------------------------------------------------------------------
# Block A: new paragraph inserted right before "TASK SCHEDULER Tips:" -
# a numbered ListParagraph with the "Write a script ... substantial." text,
# including the grammar-checker proofErr bookends around "Usually".
# ---------------------------------------------------------------------------
$blockA = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:line="276" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:rPr><w:u w:val="double"/></w:rPr><w:t>Write a script</w:t></w:r><w:r><w:t xml:space="preserve"> for 1 or more tokens to be monitored for sudden change; ex, those which have climbed in the top 10 and have potential increase in the next hours. Monitor mainly for candlestick increase with volume simultaneously. Add also MA. Give weights to indicators. This code is very IMP since new comers always experience several spikes per day before leaving the top 10. The script is critical to catch those moments again. </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Usually</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> their gains are substantial.</w:t></w:r></w:p>
'@

# Block B: a blank ListParagraph (no numbering), mirrors the blank line that
# already precedes the TASK SCHEDULER paragraph.
$blockB = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p>
'@

$insertBefore = $d.Range($taskStart, $taskStart)
$insertBefore.InsertXML($blockA + $blockB)

# The TASK SCHEDULER paragraph shifted right by the length of the inserted
# text; re-find its end so block C/D land right after it.
$findR2 = $d.Content
$found2 = $findR2.Find.Execute("TASK SCHEDULER Tips:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not re-locate 'TASK SCHEDULER Tips:' anchor text"
}
$taskPara2 = $findR2.Paragraphs.Item(1)
$taskEnd2 = $taskPara2.Range.End

# ---------------------------------------------------------------------------
# Block C: blank ListParagraph after the TASK SCHEDULER paragraph.
# Block D: blank ListParagraph with spacing/jc but no numbering (trailing
# paragraph right before the sectPr).
# ---------------------------------------------------------------------------
$blockC = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p>
'@

$blockD = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:spacing w:after="0" w:line="276" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr></w:p>
'@

$insertAfter = $d.Range($taskEnd2, $taskEnd2)
$insertAfter.InsertXML($blockC + $blockD)

Write-Output "done"
